# Scheduled price-refresh update for the Leviathan Profits workbook.
# Pulls the latest currentAveragePrice(NQ/HQ) + profit figures into
# columns H:N of each job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1111.6
$ws.Range("I2").Value = 734.875
$ws.Range("J2").Value = 1362.75
$ws.Range("K2").Value = 734.875
$ws.Range("L2").Value = 1362.75
$ws.Range("M2").Value = -621.875
$ws.Range("N2").Value = -1588.75
$ws.Range("H19").Value = 3748.8667
$ws.Range("I19").Value = 4748.4546
$ws.Range("K19").Value = 4748.4546
$ws.Range("M19").Value = -4573.4546
$ws.Range("H28").Value = 213.33333
$ws.Range("I28").Value = 213.33333
$ws.Range("K28").Value = 213.33333
$ws.Range("M28").Value = 271.66667
$ws.Range("H55").Value = 261.57895
$ws.Range("I55").Value = 236.1875
$ws.Range("K55").Value = 236.1875
$ws.Range("M55").Value = -22.1875
$ws.Range("H106").Value = 8882.8125
$ws.Range("I106").Value = 2011.8334
$ws.Range("K106").Value = 2011.8334
$ws.Range("M106").Value = -1380.8334
$ws.Range("H116").Value = 23979.8
$ws.Range("H132").Value = 1188.3208
$ws.Range("I132").Value = 1223.1459
$ws.Range("K132").Value = 3669.4377
$ws.Range("M132").Value = -1139.4377
$ws.Range("H135").Value = 29353.314
$ws.Range("I135").Value = 416.6207
$ws.Range("K135").Value = 3749.5863
$ws.Range("M135").Value = -1214.5863
$ws.Range("H137").Value = 2204.9524
$ws.Range("I137").Value = 2253.1765
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 6759.529500000001
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -4209.529500000001
$ws.Range("N137").Value = -11100

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1806.1818
$ws.Range("I2").Value = 1624.6666
$ws.Range("J2").Value = 2623
$ws.Range("K2").Value = 1624.6666
$ws.Range("L2").Value = 2623
$ws.Range("M2").Value = -1511.6666
$ws.Range("N2").Value = -2849
$ws.Range("H32").Value = 16767.771
$ws.Range("I32").Value = 3049.774
$ws.Range("K32").Value = 3049.774
$ws.Range("M32").Value = -2762.774
$ws.Range("H45").Value = 441434.53
$ws.Range("I45").Value = 632998.06
$ws.Range("K45").Value = 632998.06
$ws.Range("M45").Value = -632621.06
$ws.Range("H56").Value = 25763.334
$ws.Range("J56").Value = 14895
$ws.Range("L56").Value = 14895
$ws.Range("N56").Value = -16379
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = ""   # was -34275
$ws.Range("H110").Value = 1508.1818
$ws.Range("I110").Value = 1508.1818
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1508.1818
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 536.8181999999999
$ws.Range("N110").Value = ""   # was -6201
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").Value = ""   # was -43824
$ws.Range("H116").Value = 1806.1818
$ws.Range("I116").Value = 1624.6666
$ws.Range("J116").Value = 2623
$ws.Range("K116").Value = 1624.6666
$ws.Range("L116").Value = 2623
$ws.Range("M116").Value = 669.3334
$ws.Range("N116").Value = -7211
$ws.Range("H139").Value = 187995
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""   # was -107280

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1806.1818
$ws.Range("I3").Value = 1624.6666
$ws.Range("J3").Value = 2623
$ws.Range("K3").Value = 1624.6666
$ws.Range("L3").Value = 2623
$ws.Range("M3").Value = -1510.6666
$ws.Range("N3").Value = -2851
$ws.Range("H20").Value = 14924.056
$ws.Range("I20").Value = 15805.272
$ws.Range("J20").Value = 13539.286
$ws.Range("K20").Value = 15805.272
$ws.Range("L20").Value = 13539.286
$ws.Range("M20").Value = -15558.272
$ws.Range("N20").Value = -14033.286
$ws.Range("H86").Value = 1820.2368
$ws.Range("I86").Value = 1806.7778
$ws.Range("J86").Value = 1853.2727
$ws.Range("K86").Value = 1806.7778
$ws.Range("L86").Value = 1853.2727
$ws.Range("M86").Value = -683.7778000000001
$ws.Range("N86").Value = -4099.2727
$ws.Range("H89").Value = 1820.2368
$ws.Range("I89").Value = 1806.7778
$ws.Range("J89").Value = 1853.2727
$ws.Range("K89").Value = 9033.889000000001
$ws.Range("L89").Value = 9266.3635
$ws.Range("M89").Value = -3417.889000000001
$ws.Range("N89").Value = -20498.3635
$ws.Range("H105").Value = 3829.7778
$ws.Range("I105").Value = 4183.5
$ws.Range("K105").Value = 4183.5
$ws.Range("M105").Value = -2436.5
$ws.Range("H107").Value = 23210.596
$ws.Range("I107").Value = 31783.94
$ws.Range("K107").Value = 31783.94
$ws.Range("M107").Value = -29863.94

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 498.11765
$ws.Range("I22").Value = 397.45456
$ws.Range("K22").Value = 397.45456
$ws.Range("M22").Value = -47.45456000000001
$ws.Range("H31").Value = 55076.81
$ws.Range("I31").Value = 129374
$ws.Range("J31").Value = 22055.834
$ws.Range("K31").Value = 129374
$ws.Range("L31").Value = 22055.834
$ws.Range("M31").Value = -129079
$ws.Range("N31").Value = -22645.834
$ws.Range("H34").Value = 55076.81
$ws.Range("I34").Value = 129374
$ws.Range("J34").Value = 22055.834
$ws.Range("K34").Value = 129374
$ws.Range("L34").Value = 22055.834
$ws.Range("M34").Value = -129172
$ws.Range("N34").Value = -22459.834
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""   # was -32324
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""   # was -100968
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = ""   # was -3877
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = ""   # was -19384
$ws.Range("H105").Value = 2013.1111
$ws.Range("I105").Value = 1653
$ws.Range("K105").Value = 1653
$ws.Range("M105").Value = 94
$ws.Range("H107").Value = 1220.3914
$ws.Range("I107").Value = 953.53845
$ws.Range("K107").Value = 953.53845
$ws.Range("M107").Value = 966.46155
$ws.Range("H123").Value = 72299.1
$ws.Range("J123").Value = 72299.1
$ws.Range("L123").Value = 72299.1
$ws.Range("N123").Value = -82099.1

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1191.2142
$ws.Range("I5").Value = 1058.8334
$ws.Range("J5").Value = 1290.5
$ws.Range("K5").Value = 3176.5002
$ws.Range("L5").Value = 3871.5
$ws.Range("M5").Value = -3064.5002
$ws.Range("N5").Value = -4095.5
$ws.Range("H36").Value = 450000
$ws.Range("I36").Value = 100000
$ws.Range("J36").Value = 800000
$ws.Range("K36").Value = 300000
$ws.Range("L36").Value = 2400000
$ws.Range("M36").Value = -299831
$ws.Range("N36").Value = -2400338
$ws.Range("H54").Value = 6746.6665
$ws.Range("J54").Value = 6746.6665
$ws.Range("L54").Value = 20239.9995
$ws.Range("N54").Value = -21357.9995
$ws.Range("H135").Value = 1191.2142
$ws.Range("I135").Value = 1058.8334
$ws.Range("J135").Value = 1290.5
$ws.Range("K135").Value = 9529.5006
$ws.Range("L135").Value = 11614.5
$ws.Range("M135").Value = -6994.500599999999
$ws.Range("N135").Value = -16684.5

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 276.15625
$ws.Range("I2").Value = 266.57144
$ws.Range("K2").Value = 266.57144
$ws.Range("M2").Value = -153.57144
$ws.Range("H96").Value = 90261
$ws.Range("J96").Value = 90261
$ws.Range("L96").Value = 90261
$ws.Range("N96").Value = -95753
$ws.Range("H107").Value = 5318.1875
$ws.Range("I107").Value = 1905.25
$ws.Range("K107").Value = 1905.25
$ws.Range("M107").Value = 14.75
$ws.Range("H113").Value = 3816.6365
$ws.Range("J113").Value = 5994.5
$ws.Range("L113").Value = 5994.5
$ws.Range("N113").Value = -10334.5
$ws.Range("H126").Value = 3074.1538
$ws.Range("I126").Value = 2413.7144
$ws.Range("J126").Value = 3844.6667
$ws.Range("K126").Value = 7241.1432
$ws.Range("L126").Value = 11534.0001
$ws.Range("M126").Value = -4771.1432
$ws.Range("N126").Value = -16474.0001

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 81013.695
$ws.Range("J16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("N16").Value = -1340
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = ""   # was -25957.666
$ws.Range("I46").Value = 141673.33
$ws.Range("J46").Value = 1718.6
$ws.Range("K46").Value = 141673.33
$ws.Range("L46").Value = 1718.6
$ws.Range("M46").Value = -141485.33
$ws.Range("N46").Value = -2094.6
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = ""   # was -25125.666
$ws.Range("H60").Value = 15000
$ws.Range("I60").Value = 15000
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 15000
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -14491
$ws.Range("N60").Value = ""   # was -26018
$ws.Range("H61").Value = 86027.19
$ws.Range("I61").Value = 74293
$ws.Range("K61").Value = 74293
$ws.Range("M61").Value = -74091
$ws.Range("H113").Value = 86027.19
$ws.Range("I113").Value = 74293
$ws.Range("K113").Value = 74293
$ws.Range("M113").Value = -72123
$ws.Range("H130").Value = 85976
$ws.Range("J130").Value = 84571.2
$ws.Range("L130").Value = 84571.2
$ws.Range("N130").Value = -94611.2

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 36834.75
$ws.Range("J95").Value = 36834.75
$ws.Range("L95").Value = 36834.75
$ws.Range("N95").Value = -42326.75

Write-Output "Updated $(255) cells across $(8) sheets."
